$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table A1:M86 (header in row 1) was previously sorted descending by
# column J (log2_median_ratio). Re-sort it ascending by column M
# (wilcox_p_value) instead - this reshuffles rows 2-86 while the header row
# and all other sheet content stay put.
$dataRange = $ws.Range("A2:M86")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("M2:M86"))
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Reflect the cell the user clicked on after sorting (selection moved to M7,
# and the view was scrolled right so column F is the first visible column).
$ws.Range("M7").Select() | Out-Null
